$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new extraction run (20241014-092445 -> 20241015-090744)
$ws.Name = "IClientBalance-20241015-090744-"

# The "Dt. Referencia" column (G) moves forward one day (2024-10-14 -> 2024-10-15,
# serials 45579 -> 45580) for every data row (rows 2-274).
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45580
}

# A handful of rows also had their projected/foreseen-balance figures (and the
# matching total) revised upward in this refresh.
# Row 15: Vl. Projetado (D) 0, Saldo Previsto (E) 221.22, Vl. Total (H) 221.22
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 221.22
$ws.Cells.Item(15, 8).Value = 221.22

# Row 103: Saldo Previsto (E) 55771.97, Vl. Total (H) 55771.97
$ws.Cells.Item(103, 5).Value = 55771.97
$ws.Cells.Item(103, 8).Value = 55771.97

# Row 104: Saldo Previsto (E) 5465.97, Vl. Total (H) 5465.97
$ws.Cells.Item(104, 5).Value = 5465.97
$ws.Cells.Item(104, 8).Value = 5465.97

# Row 143: Saldo Previsto (E) 4878.87, Vl. Total (H) 4878.87
$ws.Cells.Item(143, 5).Value = 4878.87
$ws.Cells.Item(143, 8).Value = 4878.87

# Row 189: Saldo Previsto (E) 38311.72, Vl. Total (H) 38311.72
$ws.Cells.Item(189, 5).Value = 38311.72
$ws.Cells.Item(189, 8).Value = 38311.72
